# "corrección de errores cuando no hay cliente asociado"
#
# Mark the three pending tasks at the bottom of Hoja1 with their
# completion status in column C (the "% done" column used throughout
# the sheet):
#   - A61 "Error en iva cuando consumidor final en ticket"           -> 100%
#   - A68 "Cuando no hay cliente asociado no setea en ventaDTO..."   -> 100%
#   - A69 "Setear corte z en el ticket"                              -> "en proceso"
#
# Also move the on-screen selection to C70 (just past the last edited
# row), matching where the author's cursor ended up after typing these
# values in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Activate()

# Row 61: task finished -> 100% (same percentage style as the other
# "done" rows, e.g. C2, C62, C63 ...)
$ws.Range("C61").Value = 1
$ws.Range("C61").NumberFormat = "0%"

# Row 68: task finished -> 100%
$ws.Range("C68").Value = 1
$ws.Range("C68").NumberFormat = "0%"

# Row 69: task still being worked on -> "en proceso" (reuses the same
# shared string already used on C29/C66)
$ws.Range("C69").Value = "en proceso"

# Leave the view scrolled/selected where the edits finished.
$ws.Range("C70").Select()
$excel.ActiveWindow.ScrollRow = 61
$excel.ActiveWindow.ScrollColumn = 1

Write-Host "done"
